$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 799.5
$ws.Range("I19").Value = 750
$ws.Range("J19").Value = 873.75
$ws.Range("K19").Value = 750
$ws.Range("L19").Value = 873.75
$ws.Range("M19").Value = -575
$ws.Range("N19").Value = -1223.75
$ws.Range("H41").Value = 1018.4
$ws.Range("I41").Value = 1586.3334
$ws.Range("J41").Value = 166.5
$ws.Range("K41").Value = 1586.3334
$ws.Range("L41").Value = 166.5
$ws.Range("M41").Value = -1146.3334
$ws.Range("N41").Value = -1046.5
$ws.Range("H98").Value = 1725.9395
$ws.Range("I98").Value = 1269.8518
$ws.Range("K98").Value = 1269.8518
$ws.Range("M98").Value = 228.1482000000001
$ws.Range("H122").Value = 1725.9395
$ws.Range("I122").Value = 1269.8518
$ws.Range("K122").Value = 3809.5554
$ws.Range("M122").Value = -1359.5554
$ws.Range("H129").Value = 1684.4286
$ws.Range("I129").Value = 693.75
$ws.Range("J129").Value = 3005.3333
$ws.Range("K129").Value = 2081.25
$ws.Range("L129").Value = 9015.999899999999
$ws.Range("M129").Value = 2918.75
$ws.Range("N129").Value = -19015.9999
$ws.Range("H132").Value = 1562.7916
$ws.Range("I132").Value = 1333.6444
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 4000.933199999999
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -1470.933199999999
$ws.Range("N132").Value = -20060
$ws.Range("H138").Value = 2660.84
$ws.Range("I138").Value = 2435.95
$ws.Range("J138").Value = 3560.4
$ws.Range("K138").Value = 7307.849999999999
$ws.Range("L138").Value = 10681.2
$ws.Range("M138").Value = -2167.849999999999
$ws.Range("N138").Value = -20961.2

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1809.3572
$ws.Range("I2").Value = 1538.68
$ws.Range("K2").Value = 1538.68
$ws.Range("M2").Value = -1425.68
$ws.Range("H32").Value = 18516.057
$ws.Range("I32").Value = 4291.6665
$ws.Range("K32").Value = 4291.6665
$ws.Range("M32").Value = -4004.6665
$ws.Range("H61").Value = 1395.24
$ws.Range("I61").Value = 1411.7084
$ws.Range("K61").Value = 1411.7084
$ws.Range("M61").Value = -1199.7084
$ws.Range("H74").Value = 1857.8422
$ws.Range("I74").Value = 1857.8422
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 1857.8422
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -983.8422
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 1857.8422
$ws.Range("I77").Value = 1857.8422
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 9289.210999999999
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -4921.210999999999
$ws.Range("N77").ClearContents()
$ws.Range("H110").Value = 1540.5714
$ws.Range("I110").Value = 1822.7858
$ws.Range("J110").Value = 976.1429000000001
$ws.Range("K110").Value = 1822.7858
$ws.Range("L110").Value = 976.1429000000001
$ws.Range("M110").Value = 222.2141999999999
$ws.Range("N110").Value = -5066.1429
$ws.Range("H116").Value = 1809.3572
$ws.Range("I116").Value = 1538.68
$ws.Range("K116").Value = 1538.68
$ws.Range("M116").Value = 755.3199999999999
$ws.Range("H136").Value = 1395.24
$ws.Range("I136").Value = 1411.7084
$ws.Range("K136").Value = 4235.1252
$ws.Range("M136").Value = -1685.1252

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1809.3572
$ws.Range("I3").Value = 1538.68
$ws.Range("K3").Value = 1538.68
$ws.Range("M3").Value = -1424.68
$ws.Range("H86").Value = 1838.8572
$ws.Range("I86").Value = 1658.9
$ws.Range("J86").Value = 2288.75
$ws.Range("K86").Value = 1658.9
$ws.Range("L86").Value = 2288.75
$ws.Range("M86").Value = -535.9000000000001
$ws.Range("N86").Value = -4534.75
$ws.Range("H89").Value = 1838.8572
$ws.Range("I89").Value = 1658.9
$ws.Range("J89").Value = 2288.75
$ws.Range("K89").Value = 8294.5
$ws.Range("L89").Value = 11443.75
$ws.Range("M89").Value = -2678.5
$ws.Range("N89").Value = -22675.75
$ws.Range("H134").Value = 745.53845
$ws.Range("I134").Value = 754.36365
$ws.Range("J134").Value = 697
$ws.Range("K134").Value = 2263.09095
$ws.Range("L134").Value = 2091
$ws.Range("M134").Value = 271.9090500000002
$ws.Range("N134").Value = -7161

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2163.1667
$ws.Range("I16").Value = 2084
$ws.Range("J16").Value = 2559
$ws.Range("K16").Value = 2084
$ws.Range("L16").Value = 2559
$ws.Range("M16").Value = -1797
$ws.Range("N16").Value = -3133
$ws.Range("H113").Value = 2163.1667
$ws.Range("I113").Value = 2084
$ws.Range("J113").Value = 2559
$ws.Range("K113").Value = 2084
$ws.Range("L113").Value = 2559
$ws.Range("M113").Value = 86
$ws.Range("N113").Value = -6899

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 100.22222
$ws.Range("I12").Value = 8
$ws.Range("J12").Value = 126.57143
$ws.Range("K12").Value = 24
$ws.Range("L12").Value = 379.71429
$ws.Range("M12").Value = 149
$ws.Range("N12").Value = -725.71429
$ws.Range("H21").Value = 4000.25
$ws.Range("I21").Value = 2000.3334
$ws.Range("K21").Value = 6001.0002
$ws.Range("M21").Value = -5828.0002
$ws.Range("H26").Value = 182.85715
$ws.Range("I26").Value = 205
$ws.Range("J26").Value = 174
$ws.Range("K26").Value = 615
$ws.Range("L26").Value = 522
$ws.Range("M26").Value = -327
$ws.Range("N26").Value = -1098
$ws.Range("H38").Value = 61.666668
$ws.Range("I38").Value = 48.2
$ws.Range("J38").Value = 129
$ws.Range("K38").Value = 144.6
$ws.Range("L38").Value = 387
$ws.Range("M38").Value = 202.4
$ws.Range("N38").Value = -1081

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7885.4287
$ws.Range("I70").Value = 8000
$ws.Range("K70").Value = 8000
$ws.Range("M70").Value = -7730
$ws.Range("H73").Value = 7885.4287
$ws.Range("I73").Value = 8000
$ws.Range("K73").Value = 8000
$ws.Range("M73").Value = -7064
$ws.Range("H97").Value = 24095.322
$ws.Range("I97").Value = 24972.371
$ws.Range("J97").Value = 415
$ws.Range("K97").Value = 24972.371
$ws.Range("L97").Value = 415
$ws.Range("M97").Value = -24476.371
$ws.Range("N97").Value = -1407
$ws.Range("H126").Value = 2836.5264
$ws.Range("I126").Value = 2472.9333
$ws.Range("K126").Value = 7418.7999
$ws.Range("M126").Value = -4948.7999

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8453.085999999999
$ws.Range("I7").Value = 11691.4375
$ws.Range("K7").Value = 11691.4375
$ws.Range("M7").Value = -11579.4375
$ws.Range("H126").Value = 8453.085999999999
$ws.Range("I126").Value = 11691.4375
$ws.Range("K126").Value = 35074.3125
$ws.Range("M126").Value = -32604.3125

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 7496.8335
$ws.Range("I81").Value = 943
$ws.Range("J81").Value = 10773.75
$ws.Range("K81").Value = 1886
$ws.Range("L81").Value = 21547.5
$ws.Range("M81").Value = -825
$ws.Range("N81").Value = -23669.5
$ws.Range("H84").Value = 7496.8335
$ws.Range("I84").Value = 943
$ws.Range("J84").Value = 10773.75
$ws.Range("K84").Value = 9430
$ws.Range("L84").Value = 107737.5
$ws.Range("M84").Value = -4126
$ws.Range("N84").Value = -118345.5
$ws.Range("H126").Value = 1797.4333
$ws.Range("I126").Value = 1608.52
$ws.Range("J126").Value = 2742
$ws.Range("K126").Value = 4825.559999999999
$ws.Range("L126").Value = 8226
$ws.Range("M126").Value = -2355.559999999999
$ws.Range("N126").Value = -13166
$ws.Range("H132").Value = 4166.4473
$ws.Range("I132").Value = 4925.4136
$ws.Range("J132").Value = 1720.8889
$ws.Range("K132").Value = 14776.2408
$ws.Range("L132").Value = 5162.6667
$ws.Range("M132").Value = -12246.2408
$ws.Range("N132").Value = -10222.6667
